$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2857142857142857
$ws.Range("C2").Value = 0.4285714285714285
$ws.Range("P2").Value = 0.1428571428571428
$ws.Range("S2").Value = 0.1428571428571428
$ws.Range("P3").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("F6").Value = 0.2222222222222222
$ws.Range("J6").Value = 0.1111111111111111
$ws.Range("O6").Value = 0.1111111111111111
$ws.Range("R6").Value = 0.1111111111111111
$ws.Range("S6").Value = 0.4444444444444444
$ws.Range("B7").Value = 0.3333333333333333
$ws.Range("J7").Value = 0.3333333333333333
$ws.Range("Q7").Value = 0.3333333333333333
$ws.Range("F8").Value = 0.1
$ws.Range("J8").Value = 0.15
$ws.Range("O8").Value = 0.05
$ws.Range("Q8").Value = 0.1
$ws.Range("R8").Value = 0.05
$ws.Range("S8").Value = 0.55
$ws.Range("O9").Value = 0.1
$ws.Range("R9").Value = 0.2
$ws.Range("S9").Value = 0.7
$ws.Range("B10").Value = 0.08571428571428572
$ws.Range("D10").Value = 0.02857142857142857
$ws.Range("F10").Value = 0.05714285714285714
$ws.Range("J10").Value = 0.05714285714285714
$ws.Range("O10").Value = 0.05714285714285714
$ws.Range("Q10").Value = 0.1428571428571428
$ws.Range("R10").Value = 0.08571428571428572
$ws.Range("S10").Value = 0.4857142857142857
$ws.Range("K11").Value = 0.3333333333333333
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.75
$ws.Range("J12").Value = 0.25
$ws.Range("F15").Value = 0.1
$ws.Range("I15").Value = 0.1
$ws.Range("J15").Value = 0.2
$ws.Range("K15").Value = 0.2
$ws.Range("O15").Value = 0.2
$ws.Range("S15").Value = 0.2
$ws.Range("H16").Value = 0.2
$ws.Range("I16").Value = 0.2
$ws.Range("J16").Value = 0.4
$ws.Range("K16").Value = 0.2
$ws.Range("H17").Value = 0.375
$ws.Range("I17").Value = 0.5
$ws.Range("S17").Value = 0.125
$ws.Range("H18").Value = 0.2857142857142857
$ws.Range("J18").Value = 0.2857142857142857
$ws.Range("K18").Value = 0.1428571428571428
$ws.Range("S18").Value = 0.2857142857142857
$ws.Range("F19").Value = 0.01818181818181818
$ws.Range("H19").Value = 0.2727272727272727
$ws.Range("I19").Value = 0.07272727272727272
$ws.Range("J19").Value = 0.3818181818181818
$ws.Range("K19").Value = 0.01818181818181818
$ws.Range("M19").Value = 0.01818181818181818
$ws.Range("O19").Value = 0.03636363636363636
$ws.Range("S19").Value = 0.1818181818181818
